# The "No Conformidades" report tracks issues with a STATUS column (F) driven
# by a data-validation dropdown ("En proceso,Cerrada,Cancelada,Rechazada").
# This change marks every remaining "En proceso" item as "Cerrada" (closed) -
# row 10 was already "Cerrada". It also leaves the selection on F14, matching
# where the user ended up after making the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value  = "Cerrada"
$ws.Range("F7").Value  = "Cerrada"
$ws.Range("F8").Value  = "Cerrada"
$ws.Range("F9").Value  = "Cerrada"
$ws.Range("F10").Value = "Cerrada"
$ws.Range("F11").Value = "Cerrada"

$ws.Range("F14").Select()
